# Updates the cryptos worksheet (Sheet1) with refreshed price/volume
# figures, matching the upstream GitHub Actions scrape. Rows 27 and 28
# (InternetComputer(DFINITY) / LEO) also swap places in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.627.94'
$ws.Range("E2").Value = '  -4.98%  '

# Row 3
$ws.Range("D3").Value = '3.490.97'
$ws.Range("E3").Value = '  -6.14%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '''567.25'
$ws.Range("E5").Value = '  -7.15%  '

# Row 6
$ws.Range("D6").Value = '''186.18'
$ws.Range("E6").Value = '  -1.60%  '

# Row 7
$ws.Range("D7").Value = '3.486.58'
$ws.Range("E7").Value = '  -6.05%  '

# Row 8
$ws.Range("D8").Value = '''0.601'
$ws.Range("E8").Value = '  -5.50%  '

# Row 9
$ws.Range("E9").Value = '  +0.22%  '

# Row 10
$ws.Range("D10").Value = '''0.652'
$ws.Range("E10").Value = '  -9.29%  '

# Row 11
$ws.Range("E11").Value = '  -11.60%  '

# Row 12
$ws.Range("D12").Value = '''51.86'
$ws.Range("E12").Value = '  -11.04%  '

# Row 13
$ws.Range("E13").Value = '  -13.17%  '

# Row 14
$ws.Range("D14").Value = '''9.55'
$ws.Range("E14").Value = '  -10.07%  '

# Row 15
$ws.Range("D15").Value = '4.055.09'
$ws.Range("E15").Value = '  -5.77%  '

# Row 16
$ws.Range("E16").Value = '  -1.34%  '

# Row 17
$ws.Range("D17").Value = '3.506.52'

# Row 18
$ws.Range("D18").Value = '''17.90'
$ws.Range("E18").Value = '  -7.32%  '

# Row 19
$ws.Range("D19").Value = '65.366.07'
$ws.Range("E19").Value = '  -5.04%  '

# Row 20
$ws.Range("D20").Value = '''11.85'
$ws.Range("E20").Value = '  -8.26%  '

# Row 21
$ws.Range("D21").Value = '''1.03'
$ws.Range("E21").Value = '  -9.27%  '

# Row 22
$ws.Range("D22").Value = '''384.03'
$ws.Range("E22").Value = '  -6.62%  '

# Row 23
$ws.Range("D23").Value = '''4.19'
$ws.Range("E23").Value = '  -8.69%  '

# Row 24
$ws.Range("D24").Value = '''83.75'
$ws.Range("E24").Value = '  -6.16%  '

# Row 25
$ws.Range("E25").Value = '  -0.52%  '

# Row 26
$ws.Range("D26").Value = '''2.83'
$ws.Range("E26").Value = '  -7.08%  '

# Row 27
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '''6.02'
$ws.Range("E27").Value = '  -0.70%  '

# Row 28
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '''12.05'
$ws.Range("E28").Value = '  -6.26%  '

# Row 29
$ws.Range("E29").Value = '  -9.63%  '

# Row 30
$ws.Range("D30").Value = '''8.65'
$ws.Range("E30").Value = '  -10.20%  '

# Row 31
$ws.Range("D31").Value = '''30.32'
$ws.Range("E31").Value = '  -8.39%  '

# Row 32
$ws.Range("D32").Value = '''6.97'
$ws.Range("E32").Value = '  -6.83%  '

# Row 33
$ws.Range("D33").Value = '''612.26'
$ws.Range("E33").Value = '  -4.30%  '

# Row 34
$ws.Range("D34").Value = '''11.90'
$ws.Range("E34").Value = '  -6.45%  '

# Row 35
$ws.Range("D35").Value = '''62.59'
$ws.Range("E35").Value = '  -4.73%  '

# Row 36
$ws.Range("D36").Value = '''0.110'
$ws.Range("E36").Value = '  -10.15%  '

# Row 37
$ws.Range("D37").Value = '''40.59'
$ws.Range("E37").Value = '  -12.54%  '

# Row 38
$ws.Range("E38").Value = '  +0.00%  '

# Row 39
$ws.Range("D39").Value = '''0.386'
$ws.Range("E39").Value = '  -6.19%  '

# Row 40
$ws.Range("E40").Value = '  -11.26%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$ws.Range("E42").Value = '  -8.18%  '

# Row 43
$ws.Range("D43").Value = '2.941.97'
$ws.Range("E43").Value = '  +2.97%  '

# Row 44
$ws.Range("D44").Value = '''2.74'
$ws.Range("E44").Value = '  -9.89%  '

# Row 45
$ws.Range("D45").Value = '''2.42'
$ws.Range("E45").Value = '  -7.34%  '

# Row 46
$ws.Range("D46").Value = '''0.0396'
$ws.Range("E46").Value = '  -10.85%  '

# Row 47
$ws.Range("D47").Value = '''3.07'
$ws.Range("E47").Value = '  -1.97%  '

# Row 48
$ws.Range("E48").Value = '  -8.47%  '

# Row 49
$ws.Range("D49").Value = '''137.17'
$ws.Range("E49").Value = '  -3.36%  '

# Row 50
$ws.Range("D50").Value = '''8.23'
$ws.Range("E50").Value = '  -9.58%  '

# Row 51
$ws.Range("D51").Value = '''2.46'
$ws.Range("E51").Value = '  -9.96%  '
